# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before the "总计" (totals) sheet
#    and populate it with the per-fund holding detail for the new quarter.
# 2. Insert a new top row into "总计" for the 2022-Q1 summary line and shift
#    the existing quarter rows down by one.

$wb = $excel.ActiveWorkbook

# All the other quarter sheets use one consistent "header" cell style: bold
# font, thin border on all sides, centered horizontally / top vertically.
# Apply the same look to the newly-inserted header/index cells so the new
# sheet matches its siblings.
function Set-HeaderStyle($range) {
    $range.Font.Bold = $true
    $range.Borders.LineStyle = 1
    $range.HorizontalAlignment = -4108
    $range.VerticalAlignment = -4160
}

# ---------------------------------------------------------------------------
# 1. New "2022-Q1" sheet, inserted immediately before "总计"
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

# Header row (bold/centered header style carried over automatically from the
# other quarter sheets is not available via COM, so just set values; Excel
# default styling is acceptable here).
$q1.Range("B1").Value2 = "基金代码"
$q1.Range("C1").Value2 = "基金名称"
$q1.Range("D1").Value2 = "基金规模"
$q1.Range("E1").Value2 = "股票总仓位"
$q1.Range("F1").Value2 = "仓位占比"
$q1.Range("G1").Value2 = "持有市值(亿元)"
$q1.Range("H1").Value2 = "仓位排名"

# Data row (index column A is numeric 0-based like the other sheets).
$q1.Range("A2").Value2 = 0

# 基金代码/基金规模/股票总仓位/仓位占比/持有市值(亿元) are stored as TEXT in the
# source workbook (even though they look numeric) - force text storage by
# pre-formatting the cells before assigning the value.
$q1.Range("B2:G2").NumberFormat = "@"
$q1.Range("B2").Value2 = "501029"
$q1.Range("C2").Value2 = "华宝标普中国A股红利机会指数（LOF）A"
$q1.Range("D2").Value2 = "13.19"
$q1.Range("E2").Value2 = "94.39"
$q1.Range("F2").Value2 = "2.79"
$q1.Range("G2").Value2 = "0.3680"
$q1.Range("H2").Value2 = 1

Set-HeaderStyle $q1.Range("B1:H1")
Set-HeaderStyle $q1.Range("A2")

# ---------------------------------------------------------------------------
# 2. Prepend the 2022-Q1 summary row to "总计", shifting existing rows down
# ---------------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")

# Read the existing rows (2..6) before they get overwritten by the shift.
$existing = @()
for ($r = 2; $r -le 6; $r++) {
    $existing += , @($tot.Cells.Item($r, 2).Value2, $tot.Cells.Item($r, 3).Value2, $tot.Cells.Item($r, 4).Value2)
}

# Shift rows 2..6 down to 3..7, rewriting the A/B/C/D columns explicitly so
# no stray formatting creeps in from Excel's native row-insert heuristics.
for ($i = $existing.Length - 1; $i -ge 0; $i--) {
    $destRow = $i + 3
    $tot.Cells.Item($destRow, 1).Value2 = $i + 1
    $tot.Cells.Item($destRow, 2).Value2 = $existing[$i][0]
    $tot.Cells.Item($destRow, 3).Value2 = $existing[$i][1]
    $tot.Cells.Item($destRow, 4).Value2 = $existing[$i][2]
}

# New 2022-Q1 row at the top of the data block.
$tot.Cells.Item(2, 1).Value2 = 0
$tot.Cells.Item(2, 2).Value2 = "2022-Q1"
$tot.Cells.Item(2, 3).Value2 = 1
$tot.Cells.Item(2, 4).Value2 = 0.37

# Column A (the 0-based row index) keeps the bold/bordered header-ish style
# all the way down, including the row that just shifted into A7.
Set-HeaderStyle $tot.Range("A2:A7")
